$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.636.06"
$ws.Range("D3").Value = "1.643.44"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.872.53"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.673.62"
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.72%  "
$ws.Range("D17").Value = "26.687.65"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +10.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  +4.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0516"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.269.41"
$ws.Range("E35").Value = "  +4.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0180"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  +5.97%  "
$ws.Range("E39").Value = "  +2.79%  "
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.812"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").Value = "1.782.94"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.55%  "
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0974"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
